$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: reorder elements in the Scalpel Type set (content only, no style change)
$ws.Range("E2").Value = "{'any', 'plistlib.readPlistFromString'}"

# Row 3: Scalpel Type becomes 'any' (previously 'plistlib.readPlistFromString'),
# and Status flips from Win (green) to Neutral (orange) since it no longer matches.
$ws.Range("E3").Value = "any"
$ws.Range("F3").Value = "Neutral"
$ws.Range("F3").Interior.Color = 42495

# Row 4: reorder elements in the Scalpel Type set (content only, no style change)
$ws.Range("E4").Value = "{'bool', 'any'}"

# Row 5: Scalpel Type becomes 'bool' (previously 'any'),
# and Status flips from Neutral (orange) to Win (green) since it now matches.
$ws.Range("E5").Value = "bool"
$ws.Range("F5").Value = "Win"
$ws.Range("F5").Interior.Color = 32768

# Row 20: move the "Scalpel Accuracy:" label and its value two columns to the
# right so the label lines up with "Accuracy vs PyType" below it.
$ws.Range("C20").ClearContents()
$ws.Range("D20").ClearContents()
$ws.Range("E20").Value = "Scalpel Accuracy:"
$ws.Range("F20").Value = 100

# Row 21: fix wording of the accuracy label.
$ws.Range("E21").Value = "Accuracy vs PyType"
